# Fix up "Recorded By" (column G) entries on the "Session Analysis Results"
# sheet: swap "Miss Dina Nasr, Administrator" -> "Administrator, Miss Dina Nasr"
# wherever it appears, leaving plain "Miss Dina Nasr" entries untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.UsedRange.Rows.Count

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
